# Fruta / hortaliza, semanal
#
# Insert two new weekly price rows for "Damasco" / "Modesto" (Especial and
# Primera grades, fecha 2022-01-?? serial 44578) at the top of the existing
# block of rows that share that product grouping. The previously-existing
# rows 128-132 shift down to 130-134 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 128-132 down to 130-134 by inserting 2 new rows
# above the current row 128.
$ws.Rows("128:129").Insert()

# New row 128: Damasco, Modesto, Especial
$ws.Range("A128").Value = 6
$ws.Range("B128").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C128").Value = "Metropolitana"
$ws.Range("D128").Value = 44578
$ws.Range("E128").Value = 13
$ws.Range("F128").Value = "Fruta"
$ws.Range("G128").Value = 100103
$ws.Range("H128").Value = "Frutos de hueso (carozo)"
$ws.Range("I128").Value = 100103003
$ws.Range("J128").Value = "Damasco"
$ws.Range("K128").Value = "Modesto"
$ws.Range("L128").Value = "Especial"
$ws.Range("M128").Value = 170
$ws.Range("N128").Value = 17000
$ws.Range("O128").Value = 17000
$ws.Range("P128").Value = 17000
$ws.Range("Q128").Value = "`$/caja 18 kilos"
$ws.Range("R128").Value = "Región Metropolitana"
$ws.Range("S128").Value = 944
$ws.Range("T128").Value = 18

# New row 129: Damasco, Modesto, Primera
$ws.Range("A129").Value = 6
$ws.Range("B129").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C129").Value = "Metropolitana"
$ws.Range("D129").Value = 44578
$ws.Range("E129").Value = 13
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100103
$ws.Range("H129").Value = "Frutos de hueso (carozo)"
$ws.Range("I129").Value = 100103003
$ws.Range("J129").Value = "Damasco"
$ws.Range("K129").Value = "Modesto"
$ws.Range("L129").Value = "Primera"
$ws.Range("M129").Value = 250
$ws.Range("N129").Value = 14000
$ws.Range("O129").Value = 14000
$ws.Range("P129").Value = 14000
$ws.Range("Q129").Value = "`$/caja 18 kilos"
$ws.Range("R129").Value = "Región Metropolitana"
$ws.Range("S129").Value = 778
$ws.Range("T129").Value = 18
